$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Output")

# --- Update existing rows with new elastic-run results ---

# Row 2 (tech=HH2_BIO_SR_C_NEW)
$ws.Range("E2").Value = 5.392656000000001
$ws.Range("F2").Value = 34.05888000000001

# Row 3 (tech=HH2_NGA_CL_CCS_NEW)
$ws.Range("K3").Value = 818.2290596465845
$ws.Range("L3").Value = 1387.983876859377

# Row 4 (tech=HH2_COA_CL_CCS_NEW)
$ws.Range("I4").Value = 16.88011860661748
$ws.Range("J4").Value = 0.2967369502775519
$ws.Range("L4").Value = 0.354482490738961

# --- Add new row 5 for the new elastic-run technology ---

# Copy formatting/style from row 4 down to row 5 first
$ws.Range("A4:L4").Copy($ws.Range("A5:L5"))

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "HH2_COA_CM_CCS_NEW"
$ws.Range("C5").Value = "HH2_CU"
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 15.27001012745114
$ws.Range("J5").Value = 15.27001012745114
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0.3206702126764738
